$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> FAPs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 63.211268
$ws.Range("H2").Value = 189.633804
$ws.Range("I2").Value = 0.4922609885657722
$ws.Range("J2").Value = 0.4922609885657722
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7527293333333334
$ws.Range("N2").Value = 2.258188
$ws.Range("Q2").Value = 47.58097562079467
$ws.Range("R2").Value = 428.228780587152
$ws.Range("S2").Value = 0.4922609885657722
$ws.Range("T2").Value = 0.4922609885657722

# Row 3 (FAPs -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 43.30706799999999
$ws.Range("H3").Value = 129.921204
$ws.Range("I3").Value = 0.3372560111523963
$ws.Range("J3").Value = 0.3372560111523963
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.7527293333333334
$ws.Range("N3").Value = 2.258188
$ws.Range("Q3").Value = 32.59850042426133
$ws.Range("R3").Value = 293.386503818352
$ws.Range("S3").Value = 0.3372560111523963
$ws.Range("T3").Value = 0.3372560111523963

# Row 4 (sCs -> FAPs)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.891734
$ws.Range("H4").Value = 65.675202
$ws.Range("I4").Value = 0.1704830002818315
$ws.Range("J4").Value = 0.1704830002818315
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.7527293333333334
$ws.Range("N4").Value = 2.258188
$ws.Range("Q4").Value = 16.47855033933067
$ws.Range("R4").Value = 148.306953053976
$ws.Range("S4").Value = 0.1704830002818315
$ws.Range("T4").Value = 0.1704830002818315
